$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated confusion-matrix cell values (added comment for function)
$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 9
$ws.Range("F3").Value = 0
$ws.Range("S3").Value = 1
$ws.Range("X3").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 8
$ws.Range("C5").Value = 0
$ws.Range("E5").Value = 8
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 0
$ws.Range("X5").Value = 1
$ws.Range("C6").Value = 1
$ws.Range("F6").Value = 8
$ws.Range("J6").Value = 1
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 9
$ws.Range("L7").Value = 0
$ws.Range("U7").Value = 1
$ws.Range("C8").Value = 1
$ws.Range("H8").Value = 6
$ws.Range("L8").Value = 1
$ws.Range("T8").Value = 0
$ws.Range("C9").Value = 1
$ws.Range("I9").Value = 7
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 1
$ws.Range("S9").Value = 0
$ws.Range("J10").Value = 10
$ws.Range("K10").Value = 0
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 9
$ws.Range("L12").Value = 11
$ws.Range("Q12").Value = 0
$ws.Range("S12").Value = 0
$ws.Range("X12").Value = 0
$ws.Range("Y12").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("M13").Value = 10
$ws.Range("T13").Value = 0
$ws.Range("B14").Value = 0
$ws.Range("N14").Value = 10
$ws.Range("S14").Value = 0
$ws.Range("X14").Value = 1
$ws.Range("C15").Value = 1
$ws.Range("N15").Value = 3
$ws.Range("O15").Value = 6
$ws.Range("S15").Value = 0
$ws.Range("Y15").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("H16").Value = 2
$ws.Range("X16").Value = 1
$ws.Range("G17").Value = 1
$ws.Range("Q17").Value = 9
$ws.Range("S17").Value = 1
$ws.Range("H18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("R18").Value = 10
$ws.Range("X18").Value = 1
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("I19").Value = 1
$ws.Range("S19").Value = 8
$ws.Range("H20").Value = 0
$ws.Range("T20").Value = 10
$ws.Range("G21").Value = 1
$ws.Range("J21").Value = 1
$ws.Range("U21").Value = 9
$ws.Range("N22").Value = 1
$ws.Range("R22").Value = 1
$ws.Range("V22").Value = 8
$ws.Range("C23").Value = 1
$ws.Range("F23").Value = 0
$ws.Range("N23").Value = 1
$ws.Range("Q23").Value = 1
$ws.Range("W23").Value = 7
$ws.Range("N24").Value = 2
$ws.Range("S24").Value = 0
$ws.Range("X24").Value = 9
$ws.Range("B25").Value = 1
$ws.Range("L25").Value = 2
$ws.Range("Y25").Value = 8
$ws.Range("G26").Value = 1
$ws.Range("Q26").Value = 0
$ws.Range("Y26").Value = 0
$ws.Range("Z26").Value = 10
$ws.Range("B27").Value = 0
$ws.Range("J27").Value = 1
$ws.Range("AA27").Value = 9
